# [DOCS] Added API Docs link at the last slide
#
# Inserts a new "API Documentation - Postman" slide right before the
# existing closing "Thank You" slide (which stays the very last slide
# in the deck).

$p = $ppt.ActivePresentation

# The last slide in the deck is the closing "Thank You" slide. Insert
# the new slide right before it (so "Thank You" remains the final
# slide), using the "Title and Content" layout (layout #2), same as
# the rest of the content slides in this deck.
$insertAt = $p.Slides.Count
$newSlide = $p.Slides.Add($insertAt, 2)

# Title placeholder.
$titleShape = $newSlide.Shapes.Item(1)
$titleShape.TextFrame.TextRange.Text = "API Documentation - Postman"
$titleShape.TextFrame.AutoSize = 2

# Body placeholder: a single line of text that links out to the
# Postman documentation for the project's API.
$bodyShape = $newSlide.Shapes.Item(2)
$bodyTextRange = $bodyShape.TextFrame.TextRange
$bodyTextRange.Text = "Documentation Link"
$bodyTextRange.ActionSettings.Item(1).Hyperlink.Address = "https://documenter.getpostman.com/"
